$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order / data after the edit. Columns A..U.
$data = @(
    @("blooms_taxonomy", 3, 0, 3, 100, 45, 35.88, 0, 0, 3.683, 0.198, 33.33, 6.67, 100, 46.67, 0.53, 37.78, 0.0145, 17.78, 0.783, 61.93),
    @("baseline",        3, 0, 3, 100, 45, 27.6,  0, 0, 3.449, 0.161, 48.89, 15.56, 100, 42.22, 0.42, 40,    0.0101, 2.22,  0.743, 58.72),
    @("socratic",        3, 0, 3, 100, 45, 35.92, 0, 0, 5.435, 0.188, 11.11, 8.89,  100, 37.78, 0.76, 37.78, 0.0074, 15.56, 0.71,  54.89),
    @("quiz_scenario",   3, 0, 3, 100, 45, 38.12, 0, 0, 5.229, 0.167, 44.44, 26.67, 100, 42.22, 0.33, 37.78, 0.0055, 6.67,  0.81,  54.33),
    @("constrained",     3, 0, 3, 100, 45, 35.5,  0, 0, 0.378, 0.114, 13.33, 2.22,  97.78, 48.89, 0.16, 46.67, 0.0108, 2.22, 0.71,  52.59),
    @("misconceptions",  3, 2, 1, 33.33, 15, 31.02, 0, 0, 5.879, 0.173, 6.67, 0,    100, 26.67, 1.2,  20,    0.0042, 6.67,  0.519, 35.61),
    @("verbose",         3, 0, 3, 100, 45, 33.07, 0, 0, 2.856, 0.176, 11.11, 4.44, 97.78, 26.67, 0.31, 24.44, 0.0041, 22.22, 0.864, 35.18)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
